# LOQ4204.xlsx fix-up
# - Row 10 (Objetivos:) had the wrong value (the docente's name instead of
#   the actual objectives paragraph) -- replace it with the correct text.
# - A whole row was missing between "Docentes responsaveis:" (A12) and
#   "Programa resumido:" (old A13): insert it and fill in the docente's
#   name that belongs there.
# - Several B/C values were shifted/placeholder text instead of their real
#   content ("Semestral" instead of the short syllabus, a date instead of
#   the syllabus, the method text duplicated into the criteria cell, etc.)
#   -- replace each with its correct text.
# - The malformed <cols> definition (min="1" max="2" overlapping the
#   dedicated column-B width/style rule) is also normalized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objetivosPt = 'Apresentar  os conceitos básicos da Ciência Econômica, capacitando-os a compreender os principais conceitos micro e macroeconômicos e a interpretar o discurso e a prática da economia, orientados pelo seu próprio senso crítico.'
$docente = '11079086 - Herlandí de Souza Andrade'
$programaResumidoPt = 'A. Microeconomia. B. Macroeconomia. C. Desenvolvimento Econômico. D. Economia Internacional. E. Economia Brasileira'
$programaPt = 'A. MICROECONOMIA: 1. Introdução aos conceitos de Economia e fundamentos da análise microeconômica. 2. Teoria do consumidor e da demanda. 3. Teoria da firma e da oferta. 4. Custos e formação de preços. 5. Estruturas de Mercado 6. Comportamento estratégico e concorrência. 7. Tecnologia como fator de produção. 8. Sustentabilidade: recursos, custos e indicadores ambientais. B. MACROECONOMIA: 1. Fundamentos da análise macroeconômica. 2. Contabilidade nacional. 3. Equilíbrios clássicos e keynesiano. 4. Sistema monetário. 5. Política fiscal. 6. Economia mundial e comércio internacional. 7. Fundamentos da regressão como ferramenta para quantificar relações econômicas. 8. Setor público. C. DESENVOLVIMENTO ECONÔMICO: 1. Fatores de Crescimento. 2. Fontes de Desenvolvimento. 3. Financiamento do Desenvolvimento Econômico. 4. Um modelo de Crescimento Econômico. 5. O Processo de internacionalização e globalização.D. ECONOMIA INTERNACIONAL: 1. Fundamentos do Comércio Internacional. 2. Determinação das Taxas de Câmbio. 3. Políticas Externas. 4. Fatores determinantes do comportamento das importações e exportações.E. ECONOMIA BRASILEIRA: 1. A experiência histórica da industrialização brasileira. 2. A internacionalização da economia brasileira. 3. Teoria dos ciclos e realidade brasileira. 4. Os ciclos econômicos do Brasil ao longo de sua história recente.'
$bibliografiaTxt = 'VASCONCELLOS, M. A. S.; GARCIA, M. E. Fundamentos de Economia. 6 ed. São Paulo: Saraiva, 2018.GREMAUD, A. P. Introdução à Economia. São Paulo: Atlas, 2017.ROSSETTI, J. P. Introdução à Economia - Livro Texto. São Paulo: Atlas, 2016.VASCONCELLOS, M. A. S. ECONOMIA: Micro e Macro. São Paulo: Atlas, 2015.ALBERGONI, L. INTRODUÇÃO À ECONOMIA: Aplicações no Cotidiano. São Paulo: Atlas, 2015.GREMAUD, A. P.; VASCONCELLOS, M. A. S.; TONETO JÚNIOR, R. Economia Brasileira Contemporânea. 8 ed. São Paulo: Atlas, 2017.MÉNARD, C.; SAES, M. S. M.; SILVA, V. L. S.; RAYNAUD, E. Economia das Organizações: Formas Plurais e Desafios. São Paulo: Atlas, 2014.BACHA et Al. Estado da Economia Mundial - Desafios e Respostas - Seminário em Homenagem a Pedro Malan. São Paulo: LTC, 2015.BACHA , Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus,1987.BEGG, D.; DORNBUSCH, R.; FISCHER, S. Introdução A Economia. Rio de Janeiro: Campus, 2003. FURTADO, C. Formação econômica do Brasil. São Paulo: Companhia Editora Nacional, 2003.GRAMAUD, A. P. et alli. Manual de economia. São Paulo. Saraiva. 2004.GRAMAUD, A. P. et alli. Economia Brasileira Contemporânea. 6.ed. São Paulo. Atlas, 2006.HUNT, E. K.; SHERMAN, H. J. História do Pensamento Econômico. Petrópolis : Vozes, 1997.MANKIW, N.G. Introdução à economia. São Paulo: Thomson Learning, 2006.SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company.'

# 1) "Objetivos:" (row 10) was showing the docente's name -- fix it.
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# 2) Insert the missing row under "Docentes responsaveis:" (A12). This
#    shifts existing rows 13-21 down to 14-22, carrying their row heights.
$ws.Rows.Item(13).Insert()

# The freshly inserted row copied stray formatting into A13 (the
# overlapping <cols> rule bug below); the corrected sheet has no cell
# there at all, so remove it completely.
$ws.Range("A13").Clear()

# Fill in the new row with the docente's name, restoring correct column
# formatting (style 2 / style 3) by copying formats from an untouched row
# rather than re-building styles by hand.
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente
$ws.Range("B19").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) "Programa resumido:" (now row 14) incorrectly held "Semestral" --
#    replace with the real short syllabus text (Portuguese).
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt

# 4) "Programa:" (now row 16) incorrectly held a date -- replace with the
#    full Portuguese syllabus text.
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# 5) The evaluation block (now rows 19-22) had its values shifted by one:
#    "Metodo:" held the docente's name, "Criterio:" held the method text,
#    "Norma de recuperacao:" held the criteria text, and "Bibliografia:"
#    held the recovery-exam formula. Shift each value up into its correct
#    row, then fill the now-empty Bibliografia row with the real text.
$ws.Range("B19").Value = $ws.Range("B20").Value
$ws.Range("C19").Value = $ws.Range("C20").Value
$ws.Range("B20").Value = $ws.Range("B21").Value
$ws.Range("C20").Value = $ws.Range("C21").Value
$ws.Range("B21").Value = $ws.Range("B22").Value
$ws.Range("C21").Value = $ws.Range("C22").Value
$ws.Range("B22").Value = $bibliografiaTxt
$ws.Range("C22").Value = $bibliografiaTxt

# 6) Normalize the <cols> definition: column A's width/style rule should
#    only cover column A (min=1,max=1) -- it previously overlapped column
#    B's own, more specific width/style rule (min=2,max=2), which is what
#    caused new column-B cells above to inherit column A's bold style.
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(1).ColumnWidth = 30.7109375
